$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.949.31"
$ws.Range("E2").Value = "  -1.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.675.98"
$ws.Range("E3").Value = "  -2.10%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.21"
$ws.Range("E5").Value = "  -3.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.55"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  -0.45%  "

$ws.Range("E9").Value = "  -2.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -1.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.367"
$ws.Range("E11").Value = "  -3.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("E12").Value = "  -7.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.154.47"
$ws.Range("E13").Value = "  -1.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.20"
$ws.Range("E14").Value = "  -1.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.916.68"
$ws.Range("E15").Value = "  -1.04%  "

$ws.Range("E16").Value = "  -2.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.680.65"
$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.84"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.61"
$ws.Range("E19").Value = "  -3.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.52"
$ws.Range("E20").Value = "  -2.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.30"
$ws.Range("E21").Value = "  -4.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.506"
$ws.Range("E23").Value = "  -3.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.15"
$ws.Range("E24").Value = "  -1.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.17"
$ws.Range("E27").Value = "  -2.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0859"
$ws.Range("E28").Value = "  -5.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.37"
$ws.Range("E29").Value = "  +1.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.94"
$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.82"
$ws.Range("E32").Value = "  +3.02%  "

$ws.Range("E33").Value = "  +0.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.86"
$ws.Range("E34").Value = "  -0.52%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.52"
$ws.Range("E36").Value = "  -2.41%  "

$ws.Range("E37").Value = "  -0.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "350.04"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.958"
$ws.Range("E39").Value = "  -2.90%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.30"
$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.99"
$ws.Range("E41").Value = "  -2.92%  "

$ws.Range("E42").Value = "  -1.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.37"
$ws.Range("E43").Value = "  -3.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.76"
$ws.Range("E44").Value = "  -5.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0563"
$ws.Range("E45").Value = "  -3.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.617"
$ws.Range("E46").Value = "  -0.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0974"
$ws.Range("E49").Value = "  -2.75%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.07"
$ws.Range("E50").Value = "  -4.18%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0241"
$ws.Range("E51").Value = "  -2.44%  "
